# Add an additional "Sponsor Protocol Number 2" column to the upload
# template, right after the existing "Project Number 2" column (D) and
# before "Investigator Name with Qualification (ICSF)" (old E, now F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at E - shifts old E (and everything after it)
# one column to the right.
$ws.Columns("E:E").Insert() | Out-Null

# Header text for the newly inserted column.
$ws.Cells.Item(1, 5).Value = "Sponsor Protocol Number 2"

# Size the new column to fit its header text, matching the neighbouring
# bestFit-style header columns.
$ws.Columns("E:E").ColumnWidth = 24.86

# Leave the new header cell selected, as in the authored workbook.
$ws.Range("E1").Select() | Out-Null
